# Helper: set a cell's value as forced text (avoids numeric auto-conversion,
# e.g. leading zeros in fund codes, or trailing-zero decimals like "0.40"),
# then reset the cell style to the workbook default so no stray numeric
# format / quote-prefix style is left on the cell.
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计".
#    The cleanest way to get an exact match of the existing detail-sheet
#    formatting (header row + index column both use the same bold/bordered
#    style) is to duplicate the existing "2022-Q2" sheet (which has the
#    same header labels) and then overwrite its data cells.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($null, $totalSheet)

$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q3"

# The copied sheet has 10 data rows (rows 2-11); 2022-Q3 only has 9 (rows 2-10)
$new.Rows.Item(11).Delete()

Set-TextCell $new "B2" "002376"
Set-TextCell $new "C2" "国寿安保核心产业灵活配置混合"
Set-TextCell $new "D2" "4.85"
Set-TextCell $new "E2" "86.85"
Set-TextCell $new "F2" "3.63"
Set-TextCell $new "G2" "0.1761"
$new.Range("H2").Value = 7

Set-TextCell $new "B3" "005175"
Set-TextCell $new "C3" "国寿安保消费新蓝海灵活配置混合"
Set-TextCell $new "D3" "0.72"
Set-TextCell $new "E3" "90.95"
Set-TextCell $new "F3" "4.52"
Set-TextCell $new "G3" "0.0325"
$new.Range("H3").Value = 4

Set-TextCell $new "B4" "004258"
Set-TextCell $new "C4" "国寿安保稳嘉混合A"
Set-TextCell $new "D4" "2.32"
Set-TextCell $new "E4" "20.12"
Set-TextCell $new "F4" "0.99"
Set-TextCell $new "G4" "0.0230"
$new.Range("H4").Value = 5

Set-TextCell $new "B5" "004301"
Set-TextCell $new "C5" "国寿安保稳信混合A"
Set-TextCell $new "D5" "1.50"
Set-TextCell $new "E5" "22.10"
Set-TextCell $new "F5" "0.99"
Set-TextCell $new "G5" "0.0148"
$new.Range("H5").Value = 6

Set-TextCell $new "B6" "012665"
Set-TextCell $new "C6" "国寿安保裕祥混合A"
Set-TextCell $new "D6" "0.12"
Set-TextCell $new "E6" "34.24"
Set-TextCell $new "F6" "1.55"
Set-TextCell $new "G6" "0.0019"
$new.Range("H6").Value = 8

Set-TextCell $new "B7" "012666"
Set-TextCell $new "C7" "国寿安保裕祥混合C"
Set-TextCell $new "D7" "0.03"
Set-TextCell $new "E7" "34.24"
Set-TextCell $new "F7" "1.55"
Set-TextCell $new "G7" "0.0005"
$new.Range("H7").Value = 8

Set-TextCell $new "B8" "004302"
Set-TextCell $new "C8" "国寿安保稳信混合C"
Set-TextCell $new "D8" "0.02"
Set-TextCell $new "E8" "22.10"
Set-TextCell $new "F8" "0.99"
Set-TextCell $new "G8" "0.0002"
$new.Range("H8").Value = 6

Set-TextCell $new "B9" "015406"
Set-TextCell $new "C9" "国寿安保稳信混合E"
Set-TextCell $new "D9" "0.01"
Set-TextCell $new "E9" "22.10"
Set-TextCell $new "F9" "0.99"
Set-TextCell $new "G9" "0.0001"
$new.Range("H9").Value = 6

Set-TextCell $new "B10" "004259"
Set-TextCell $new "C10" "国寿安保稳嘉混合C"
Set-TextCell $new "D10" "0.00"
Set-TextCell $new "E10" "20.12"
Set-TextCell $new "F10" "0.99"
$new.Range("G10").Value = 0
$new.Range("H10").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 with the
#    2022-Q3 totals, pushing the existing quarters down by one row.
# ------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# restore the index-column style (bold/bordered) on the new A2 cell by
# copying the format already present on A3 (the row that was A2 before
# the insert, and already carries that style)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet "B2" "2022-Q3"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 0.25

# the rows that shifted down (old rows 2-6, now rows 3-7) kept their old
# 0-based index values in column A; renumber them to stay sequential
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# restore original active sheet to avoid leaving a stray tabSelected flag
$totalSheet.Activate()
